$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update vendor address / contact info for rows 9 and 10
# (written in this column order so new shared-string entries are created
# in the same sequence Excel used when the workbook was authored)
# Columns T, U, X, Y are number-formatted/quote-prefixed text cells, so a
# leading apostrophe is used to keep them stored as quote-prefixed text
# (preserving their existing cell style) instead of Excel re-styling them.
$ws.Range("X9").Formula = "'201-999-5654"
$ws.Range("Y9").Formula = "'210-339-0102"
$ws.Range("S9").Value = "Street Sacramento"
$ws.Range("U9").Formula = "'USA"
$ws.Range("T9").Formula = "'Alaska"
$ws.Range("V9").Value = 99501

$ws.Range("X10").Formula = "'201-999-5654"
$ws.Range("Y10").Formula = "'210-339-0102"
$ws.Range("S10").Value = "Street Sacramento"
$ws.Range("U10").Formula = "'USA"
$ws.Range("T10").Formula = "'Alaska"
$ws.Range("V10").Value = 99501

# Update the selected/active view state to match the saved workbook view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 18
$ws.Range("W9").Select()
